$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.3821403333333333
$ws.Range("H2").Value = 1.146421
$ws.Range("I2").Value = 0.06500081136128052
$ws.Range("J2").Value = 0.06500081136128051
$ws.Range("M2").Value = 0.05057900000000001
$ws.Range("N2").Value = 0.151737
$ws.Range("O2").Value = 0.01400296657613869
$ws.Range("P2").Value = 0.01400296657613869
$ws.Range("Q2").Value = 0.01932827591966667
$ws.Range("R2").Value = 0.173954483277
$ws.Range("S2").Value = 0.0009102041889139073
$ws.Range("T2").Value = 0.000910204188913907

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.3821403333333333
$ws.Range("H3").Value = 1.146421
$ws.Range("I3").Value = 0.06500081136128052
$ws.Range("J3").Value = 0.06500081136128051
$ws.Range("O3").Value = 0.146324388539341
$ws.Range("P3").Value = 0.146324388539341
$ws.Range("Q3").Value = 0.2019713565755556
$ws.Range("R3").Value = 1.81774220918
$ws.Range("S3").Value = 0.009511203977000421
$ws.Range("T3").Value = 0.009511203977000419

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.3821403333333333
$ws.Range("H4").Value = 1.146421
$ws.Range("I4").Value = 0.06500081136128052
$ws.Range("J4").Value = 0.06500081136128051
$ws.Range("O4").Value = 0.8396726448845202
$ws.Range("P4").Value = 0.8396726448845202
$ws.Range("Q4").Value = 1.158999021691556
$ws.Range("R4").Value = 10.430991195224
$ws.Range("S4").Value = 0.05457940319536619
$ws.Range("T4").Value = 0.05457940319536618

$ws.Range("I5").Value = 0.8014409005237051
$ws.Range("J5").Value = 0.801440900523705
$ws.Range("M5").Value = 0.05057900000000001
$ws.Range("N5").Value = 0.151737
$ws.Range("O5").Value = 0.01400296657613869
$ws.Range("P5").Value = 0.01400296657613869
$ws.Range("Q5").Value = 0.238311961562
$ws.Range("R5").Value = 2.144807654058
$ws.Range("S5").Value = 0.01122255014278394
$ws.Range("T5").Value = 0.01122255014278393

$ws.Range("I6").Value = 0.8014409005237051
$ws.Range("J6").Value = 0.801440900523705
$ws.Range("O6").Value = 0.146324388539341
$ws.Range("P6").Value = 0.146324388539341
$ws.Range("S6").Value = 0.11727034971955
$ws.Range("T6").Value = 0.1172703497195499

$ws.Range("I7").Value = 0.8014409005237051
$ws.Range("J7").Value = 0.801440900523705
$ws.Range("O7").Value = 0.8396726448845202
$ws.Range("P7").Value = 0.8396726448845202
$ws.Range("S7").Value = 0.6729480006613711
$ws.Range("T7").Value = 0.6729480006613711

$ws.Range("I8").Value = 0.1335582881150144
$ws.Range("J8").Value = 0.1335582881150144
$ws.Range("M8").Value = 0.05057900000000001
$ws.Range("N8").Value = 0.151737
$ws.Range("O8").Value = 0.01400296657613869
$ws.Range("P8").Value = 0.01400296657613869
$ws.Range("Q8").Value = 0.03971414186966667
$ws.Range("R8").Value = 0.357427276827
$ws.Range("S8").Value = 0.001870212244440848
$ws.Range("T8").Value = 0.001870212244440848

$ws.Range("I9").Value = 0.1335582881150144
$ws.Range("J9").Value = 0.1335582881150144
$ws.Range("O9").Value = 0.146324388539341
$ws.Range("P9").Value = 0.146324388539341
$ws.Range("S9").Value = 0.01954283484279062
$ws.Range("T9").Value = 0.01954283484279061

$ws.Range("I10").Value = 0.1335582881150144
$ws.Range("J10").Value = 0.1335582881150144
$ws.Range("O10").Value = 0.8396726448845202
$ws.Range("P10").Value = 0.8396726448845202
$ws.Range("S10").Value = 0.1121452410277829
$ws.Range("T10").Value = 0.1121452410277829

Write-Output "Applied TPM update changes"